$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Insert two new rows after row 4 (Adobe Creative Cloud's EdgeWebView row),
# shifting every row from the old row 5 onward down by two.
$ws.Rows.Item(5).Resize(2).Insert()

# Populate the two newly inserted rows with the "Adobe Desktop" entry.
$ws.Range("A5").Value = "##Adobe Desktop"
$ws.Range("A6").Value = "%ProgramFilesFolder32%\Common Files\Adobe\Adobe Desktop Common"
$ws.Range("E6").Value = "0x00000221"

# Column width tweaks that accompanied the new, longer path text.
$ws.Columns.Item(1).ColumnWidth = 61.5
$ws.Columns.Item(2).ColumnWidth = 27.25

# Selection moved to A12 in the saved workbook.
$ws.Range("A12").Select()
